# Excel COM-interop script: EntProcessor for graph_builder
#
# 1. Rewrite the "Ents" column (J2:J40) shared-string JSON payloads from the
#    old nested-array form  [[[...factors...], [...outcomes...]], ...]
#    to the new list-of-objects form
#    [{"factor": [...], "outcome": [...]}, ...]
# 2. Turn on AutoFilter for the used range (A1:J40) and add the resulting
#    hidden _FilterDatabase defined name.
# 3. Leave the selection on the last-touched cell (J39), matching the
#    workbook's saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = '[{"factor": ["serum", "thioredoxin concentration"], "outcome": ["imaes"]}]'
$ws.Range("J3").Value = '[{"factor": ["tbi"], "outcome": ["mortality"]}, {"factor": ["age"], "outcome": ["mortality"]}, {"factor": ["glasgow coma scale"], "outcome": ["mortality"]}, {"factor": [], "outcome": ["mortality"]}, {"factor": ["anysochoria"], "outcome": ["mortality"]}]'
$ws.Range("J4").Value = '[{"factor": ["mif"], "outcome": ["in-hospital major adverse event", "long-term", "clinical outcome", "severity", "inflammation"]}]'
$ws.Range("J5").Value = '[{"factor": ["glasgow coma scale score"], "outcome": ["gos"]}, {"factor": ["smr", "mortality rate", "standard"], "outcome": ["mortality rate"]}, {"factor": ["gos"], "outcome": ["years", "post-trauma", "gos"]}]'
$ws.Range("J6").Value = '[{"factor": ["analysis", "impact", "prognostic model", "prognosis", "international mission", "clinical trial", "tbi"], "outcome": ["severe", "prognosis", "patient", "tbi", "elderly"]}]'
$ws.Range("J7").Value = '[{"factor": ["time post-tbi"], "outcome": ["percentage", "cognitive functionality gain"]}, {"factor": ["admission", "cognitive function"], "outcome": ["percentage", "cognitive functionality gain"]}]'
$ws.Range("J8").Value = '[{"factor": ["ct"], "outcome": ["death"]}, {"factor": ["marshall and rotterdam scoring system"], "outcome": ["death"]}, {"factor": ["basal cistern absence"], "outcome": ["death"]}, {"factor": ["positive midline shift"], "outcome": ["death"]}, {"factor": ["hemorrhagic mass"], "outcome": ["death"]}, {"factor": ["subarachnoid hemorrhage", "intraventricular"], "outcome": ["death"]}]'
$ws.Range("J9").Value = '[{"factor": ["apache ii"], "outcome": ["icu-treated", "patient", "mortality", "tbi", "six-month"]}, {"factor": ["sap ii"], "outcome": ["icu-treated", "patient", "mortality", "tbi", "six-month"]}, {"factor": ["sofa"], "outcome": ["icu-treated", "patient", "mortality", "tbi", "six-month"]}, {"factor": ["age"], "outcome": ["icu-treated", "patient", "mortality", "tbi", "six-month"]}, {"factor": ["glasgow coma scale"], "outcome": ["icu-treated", "patient", "mortality", "tbi", "six-month"]}]'
$ws.Range("J10").Value = '[{"factor": ["v/c ratio"], "outcome": ["gos score", "drs", "lcf score"]}]'
$ws.Range("J11").Value = '[{"factor": ["serum", "timp-1 level"], "outcome": ["mortality"]}]'
$ws.Range("J12").Value = '[{"factor": ["evidence", "mri", "contusion"], "outcome": ["gos-e"]}, {"factor": ["roi", "reduce", "fa", "severely"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["age"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}, {"factor": ["roi", "reduce", "fa", "severely"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}]'
$ws.Range("J13").Value = '[{"factor": ["preinjury", "depressive symptom"], "outcome": ["affective", "behavioral", "problem"]}, {"factor": ["preinjury", "depressive symptom"], "outcome": ["cognitive problem"]}, {"factor": ["preinjury", "depressive symptom"], "outcome": ["mental health-related quality-of-life"]}]'
$ws.Range("J14").Value = '[{"factor": ["age"], "outcome": ["outcome", "unfavorable"]}, {"factor": ["admission", "gcs"], "outcome": []}, {"factor": ["injury", "mechanism"], "outcome": []}, {"factor": ["iss"], "outcome": []}, {"factor": ["ais", "head"], "outcome": []}, {"factor": ["type", "operation"], "outcome": []}, {"factor": ["hemorrhage acuity"], "outcome": []}, {"factor": ["operation", "time"], "outcome": []}, {"factor": ["warfarin", "clopidogrel", "pre-hospital"], "outcome": []}, {"factor": ["in-hospital death"], "outcome": []}]'
$ws.Range("J15").Value = '[{"factor": ["old age"], "outcome": ["poor"]}, {"factor": ["female gender"], "outcome": ["poor"]}]'
$ws.Range("J16").Value = '[{"factor": ["age"], "outcome": ["surgery", "outcome", "tbi"]}, {"factor": ["sex"], "outcome": ["surgery", "outcome", "tbi"]}, {"factor": ["score", "deyo-charlson comorbidity index"], "outcome": ["surgery", "outcome", "tbi"]}, {"factor": ["hospital volume"], "outcome": ["surgery", "outcome", "tbi"]}, {"factor": ["volume", "surgeon"], "outcome": ["surgery", "outcome", "tbi"]}]'
$ws.Range("J17").Value = '[{"factor": ["age"], "outcome": ["mortality", "icu"]}, {"factor": ["gender"], "outcome": ["mortality", "icu"]}, {"factor": ["etiology", "injury"], "outcome": ["mortality", "icu"]}]'
$ws.Range("J18").Value = '[{"factor": ["csf", "concentration"], "outcome": ["mortality"]}, {"factor": ["concentration", "plasma"], "outcome": ["mortality"]}, {"factor": ["abeta42", "csf", "concentration", "change"], "outcome": ["neurological status"]}]'
$ws.Range("J19").Value = '[{"factor": ["soluble", "plasminogen activator receptor", "supar", "urokinase"], "outcome": ["tbi"]}]'
$ws.Range("J20").Value = '[{"factor": ["extend", "model", "core", "impact"], "outcome": ["outcome", "month"]}, {"factor": ["basic model", "crash"], "outcome": ["outcome", "month"]}, {"factor": ["nijmegen model"], "outcome": ["outcome", "month"]}, {"factor": ["extend", "model", "core", "impact"], "outcome": ["mortality", "month"]}, {"factor": ["basic model", "crash"], "outcome": ["mortality", "month"]}, {"factor": ["nijmegen model"], "outcome": ["mortality", "month"]}]'
$ws.Range("J21").Value = '[{"factor": ["coagulopathy", "abnormal", "fibrinolysis"], "outcome": ["surgery", "deterioration", "traumatic brain injury"]}]'
$ws.Range("J22").Value = '[{"factor": ["cct", "central conduction time"], "outcome": ["long-term", "clinical outcome"]}, {"factor": ["latency"], "outcome": ["long-term", "clinical outcome"]}]'
$ws.Range("J23").Value = '[{"factor": ["sbp"], "outcome": ["mortality"]}]'
$ws.Range("J24").Value = '[{"factor": ["d20"], "outcome": ["death"]}, {"factor": ["d25"], "outcome": ["death"]}, {"factor": ["dprx"], "outcome": ["death"]}]'
$ws.Range("J25").Value = '[{"factor": ["tsp-1", "level"], "outcome": ["mortality", "1-week"]}, {"factor": ["tsp-1", "level"], "outcome": ["mortality"]}, {"factor": ["tsp-1", "level"], "outcome": ["outcome", "unfavorable"]}]'
$ws.Range("J26").Value = '[{"factor": ["bdnf", "plasma level"], "outcome": ["severe", "patient", "mortality", "tbi", "icu"]}]'
$ws.Range("J27").Value = '[{"factor": ["crash-ct model"], "outcome": ["death", "days"]}]'
$ws.Range("J28").Value = '[{"factor": ["time to death"], "outcome": ["life-sustaining", "withdrawal", "therapy"]}, {"factor": ["center", "classification"], "outcome": ["mortality"]}]'
$ws.Range("J29").Value = '[{"factor": ["prognostic model", "impact"], "outcome": ["patient", "mortality", "unfavorable outcome", "tbi"]}]'
$ws.Range("J30").Value = '[{"factor": ["serum", "concentration", "trx"], "outcome": ["mortality", "1-week"]}, {"factor": [], "outcome": ["mortality"]}, {"factor": [], "outcome": ["outcome", "unfavorable"]}]'
$ws.Range("J31").Value = '[{"factor": ["gcs"], "outcome": ["severity", "tbi"]}, {"factor": ["duration", "pta"], "outcome": ["severity", "tbi"]}, {"factor": ["assault", "motor vehicle collision", "victim of fall", "victim"], "outcome": ["parosmia score"]}]'
$ws.Range("J32").Value = '[{"factor": ["il-6", "level"], "outcome": ["septic", "development"]}, {"factor": ["il-6", "level"], "outcome": ["multiple organ dysfunction", "development"]}, {"factor": ["c-reactive protein level"], "outcome": ["septic", "development"]}, {"factor": ["c-reactive protein level"], "outcome": ["multiple organ dysfunction", "development"]}, {"factor": ["leukocyte count"], "outcome": ["septic", "development"]}, {"factor": ["leukocyte count"], "outcome": ["multiple organ dysfunction", "development"]}]'
$ws.Range("J33").Value = '[{"factor": ["rotterdam"], "outcome": ["weeks", "death"]}]'
$ws.Range("J34").Value = '[{"factor": ["admission serum glucose level"], "outcome": ["in-hospital mortality rate"]}]'
$ws.Range("J35").Value = '[{"factor": ["inflammatory mediator", "pca", "tbi"], "outcome": ["outcome"]}]'
$ws.Range("J36").Value = '[{"factor": ["age"], "outcome": ["poor outcome"]}, {"factor": ["gcs"], "outcome": ["poor outcome"]}, {"factor": ["iss"], "outcome": ["poor outcome"]}, {"factor": ["ais", "head"], "outcome": ["poor outcome"]}, {"factor": ["sex"], "outcome": ["poor outcome"]}, {"factor": ["mot"], "outcome": ["poor outcome"]}, {"factor": ["lung injury"], "outcome": ["poor outcome"]}, {"factor": ["severity", "lung injury"], "outcome": ["poor outcome"]}]'
$ws.Range("J37").Value = '[{"factor": ["il-6"], "outcome": ["year", "favorable", "gos"]}, {"factor": ["pg"], "outcome": ["year", "favorable", "gos"]}, {"factor": ["gfap"], "outcome": ["unfavorable", "year", "score", "gos"]}, {"factor": ["il-6"], "outcome": ["year", "survival status"]}, {"factor": ["pg"], "outcome": ["year", "survival status"]}, {"factor": ["gfap"], "outcome": ["year", "survival status"]}]'
$ws.Range("J38").Value = '[{"factor": ["acute"], "outcome": ["score", "drs"]}, {"factor": ["fa", "subacute"], "outcome": ["score", "drs"]}]'
$ws.Range("J39").Value = '[{"factor": ["mechanical ventilation"], "outcome": ["neurological"]}, {"factor": ["severity", "head injury"], "outcome": ["neurological"]}, {"factor": ["blood transfusion"], "outcome": ["neurological"]}, {"factor": ["neurosurgical intervention"], "outcome": ["neurological"]}, {"factor": ["mechanical ventilation"], "outcome": ["non-neurological", "complication"]}, {"factor": ["glasgow coma scale"], "outcome": ["non-neurological", "complication"]}, {"factor": ["blood transfusion"], "outcome": ["non-neurological", "complication"]}, {"factor": ["injury", "concomitant"], "outcome": ["non-neurological", "complication"]}]'
$ws.Range("J40").Value = '[{"factor": ["traumatic brain injury"], "outcome": ["mortality"]}]'

# Re-apply AutoFilter on the table range (A1:J40) -- adds <autoFilter ref="A1:J40".../>
# to the worksheet and a hidden _xlnm._FilterDatabase defined name scoped to Sheet1.
$ws.Range("A1:J40").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$J`$40")
$filterName.Visible = $false

# Move the view / selection to match the saved cursor position (row 22 scrolled
# into view, active cell J39).
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J39").Select() | Out-Null
